$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.075.81'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '1.775.24'
$ws.Range('E3').Value = '  -2.47%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.30'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.46%  '
$ws.Range('E6').Value = '  +0.25%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '31.65'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.73%  '
$ws.Range('E9').Value = '  -0.41%  '
$ws.Range('E10').Value = '  -1.86%  '
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('D12').Value = '2.031.37'
$ws.Range('E12').Value = '  -2.36%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.00'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +5.76%  '
$ws.Range('D14').Value = '1.779.56'
$ws.Range('E14').Value = '  -2.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.624'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -3.31%  '
$ws.Range('D16').Value = '34.010.19'
$ws.Range('E16').Value = '  -0.19%  '
$ws.Range('E17').Value = '  -1.95%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.75'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '253.28'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.98%  '
$ws.Range('D20').Value = '0.0₃0737'
$ws.Range('E20').Value = '  -1.90%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.999'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.34'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.05%  '
$ws.Range('E23').Value = '  -3.56%  '
$ws.Range('E24').Value = '  -2.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '156.57'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.98%  '
$ws.Range('E26').Value = '  -1.05%  '
$ws.Range('E27').Value = '  -2.82%  '
$ws.Range('E28').Value = '  -1.27%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').Value = '  -3.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0510'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.62%  '
$ws.Range('E32').Value = '  -0.69%  '
$ws.Range('E33').Value = '  +1.15%  '
$ws.Range('E34').Value = '  +2.06%  '
$ws.Range('D35').Value = '1.444.30'
$ws.Range('E35').Value = '  -6.54%  '
$ws.Range('E36').Value = '  -4.04%  '
$ws.Range('E37').Value = '  -0.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.624'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.16%  '
$ws.Range('E39').Value = '  +1.43%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '82.73'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.92%  '
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.888'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.92%  '
$ws.Range('E43').Value = '  -5.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0507'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.57%  '
$ws.Range('E45').Value = '  -1.79%  '
$ws.Range('D46').Value = '1.930.13'
$ws.Range('E46').Value = '  -2.66%  '
$ws.Range('E47').Value = '  +0.81%  '
$ws.Range('E48').Value = '  +1.71%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '98.13'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.77%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '49.45'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -6.16%  '
